# Correct destination links Update version control date
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace the outdated "National Pupil Database" permalinks in column B (rows 11 & 12)
# with the corrected destination-measures publication links.
$ws.Range("B11").Value = "<a href = 'https://explore-education-statistics.service.gov.uk/find-statistics/key-stage-4-destination-measures/2022-23'>Key stage 4 destination measures</a>"
$ws.Range("B12").Value = "<a href = 'https://explore-education-statistics.service.gov.uk/find-statistics/16-18-destination-measures/2022-23'>16-18 destination measures</a>"

# Move the active cell selection to B13 (previously C13)
$ws.Range("B13").Select()
